# issue #5: stock data output to json file
# The "股票" (stock) worksheet gains a "property_category" column (value
# "stock" for every data row) inserted right before the existing "date"
# column, and a new "legislator_id" column (828) appended at the end.
# Also fixes a stray space in a company name: "中國合成橡膠股份有限公 司"
# -> "中國合成橡膠股份有限公司".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Fix typo'd company name (extra space) before anything else.
$ws.Cells.Item(7, 2).Value = "中國合成橡膠股份有限公司"

# Insert a new column H ("property_category") ahead of the existing
# H ("date") column; everything from H onward shifts right by one.
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
